$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; Text="71÷6=11, 5"},
    @{Row=1;  Col=2; Text="18÷3=6, 0"},
    @{Row=1;  Col=3; Text="34÷9=3, 7"},
    @{Row=1;  Col=4; Text="30÷3=10, 0"},
    @{Row=1;  Col=5; Text="72÷3=24, 0"},

    @{Row=5;  Col=1; Text="59÷8=7, 3"},
    @{Row=5;  Col=2; Text="50÷9=5, 5"},
    @{Row=5;  Col=3; Text="41÷4=10, 1"},
    @{Row=5;  Col=4; Text="16÷6=2, 4"},
    @{Row=5;  Col=5; Text="10÷4=2, 2"},

    @{Row=9;  Col=1; Text="48÷7=6, 6"},
    @{Row=9;  Col=2; Text="52÷2=26, 0"},
    @{Row=9;  Col=3; Text="87÷2=43, 1"},
    @{Row=9;  Col=4; Text="15÷2=7, 1"},
    @{Row=9;  Col=5; Text="75÷8=9, 3"},

    @{Row=13; Col=1; Text="65÷5=13, 0"},
    @{Row=13; Col=2; Text="30÷2=15, 0"},
    @{Row=13; Col=3; Text="61÷3=20, 1"},
    @{Row=13; Col=4; Text="81÷4=20, 1"},
    @{Row=13; Col=5; Text="78÷8=9, 6"},

    @{Row=17; Col=1; Text="60÷3=20, 0"},
    @{Row=17; Col=2; Text="24÷7=3, 3"},
    @{Row=17; Col=3; Text="50÷3=16, 2"},
    @{Row=17; Col=4; Text="34÷3=11, 1"},
    @{Row=17; Col=5; Text="70÷6=11, 4"}
)

foreach ($r in $replacements) {
    $cell = $t.Cell($r.Row, $r.Col)
    $cell.Range.Text = $r.Text
}
